$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 1449.065440567477
    3 = 1682.753624963989
    4 = 1447.771053132694
    5 = 1152.971156719062
    6 = 920.9791591586903
    7 = 743.1989985483677
    8 = 597.3139524519362
    9 = 476.6218033999108
    10 = 375.5045118688228
    11 = 289.6988298118759
    12 = 1520.321676520881
    13 = 1758.337554657989
    14 = 1551.515582945341
    15 = 1301.993310865502
    16 = 1078.703789113911
    17 = 893.936180847612
    18 = 742.5533736725284
    19 = 617.467448569966
    20 = 512.7732686154354
    21 = 424.0048685738068
    22 = 1675.191658989845
    23 = 1494.721883091467
    24 = 1248.611381691546
    25 = 1028.201902828336
    26 = 845.7003259388883
    27 = 696.0985646797196
    28 = 572.434850711369
    29 = 468.8976854109573
    30 = 381.087104637552
    31 = 1506.109577995753
    32 = 1332.950047786653
    33 = 1096.311567089002
    34 = 883.9216920269789
    35 = 707.7453695961514
    36 = 563.124656019628
    37 = 443.4444915726271
    38 = 343.1520477882024
    39 = 258.0302704944485
    40 = 1137.361745777222
    41 = 911.703372105219
    42 = 708.6539849279947
    43 = 539.8734079048023
    44 = 401.0946340445006
    45 = 286.0979061870754
    46 = 189.6277732654353
    47 = 107.6785533116764
    48 = 941.4207343773982
    49 = 726.265051145574
    50 = 532.1945877320786
    51 = 370.5527283406426
    52 = 237.4312466326659
    53 = 126.9816923928614
    54 = 34.23044991263544
    55 = -44.6267212899925
    56 = 552.9903076125372
    57 = 366.9589104992481
    58 = 211.7330047230122
    59 = 83.71151428529014
    60 = -22.62884955591822
    61 = -112.0127849491428
    62 = -188.065682112287
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
